$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 (year 2025) metrics
$ws.Range("C8").Value = 1167
$ws.Range("D8").Value = 193
$ws.Range("E8").Value = 974
$ws.Range("F8").Value = 7.916324856439704
$ws.Range("G8").Value = 83.46186803770351
$ws.Range("H8").Value = 16.53813196229649
